# Scheduled price/profit-sheet refresh for Coeurl_Profits workbook.
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) for a set of
# rows across all eight crafter-job sheets with freshly pulled market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 119.71429
$ws.Range("I49").Value = 137.25
$ws.Range("J49").Value = 96.333336
$ws.Range("K49").Value = 411.75
$ws.Range("L49").Value = 289.000008
$ws.Range("M49").Value = -275.75
$ws.Range("N49").Value = -561.000008

$ws.Range("H69").Value = 9967.571
$ws.Range("I69").Value = 3225
$ws.Range("J69").Value = 11091.333
$ws.Range("K69").Value = 9675
$ws.Range("L69").Value = 33273.999
$ws.Range("M69").Value = -8801
$ws.Range("N69").Value = -35021.999

$ws.Range("H72").Value = 9967.571
$ws.Range("I72").Value = 3225
$ws.Range("J72").Value = 11091.333
$ws.Range("K72").Value = 29025
$ws.Range("L72").Value = 99821.997
$ws.Range("M72").Value = -24657
$ws.Range("N72").Value = -108557.997

$ws.Range("H116").Value = 9351.375
$ws.Range("I116").Value = 8773.294
$ws.Range("K116").Value = 8773.294
$ws.Range("M116").Value = -5331.294

$ws.Range("H136").Value = 82169.5
$ws.Range("J136").Value = 82169.5
$ws.Range("L136").Value = 82169.5
$ws.Range("N136").Value = -92369.5

$ws.Range("H139").Value = 137963.5
$ws.Range("J139").Value = 146945.25
$ws.Range("L139").Value = 146945.25
$ws.Range("N139").Value = -157225.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2417.8823
$ws.Range("I122").Value = 2144.2856
$ws.Range("K122").Value = 6432.8568
$ws.Range("M122").Value = -3982.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5087.9473
$ws.Range("I99").Value = 2136.3845
$ws.Range("K99").Value = 2136.3845
$ws.Range("M99").Value = -638.3845000000001

$ws.Range("H105").Value = 5219.75
$ws.Range("I105").Value = 5070.778
$ws.Range("K105").Value = 5070.778
$ws.Range("M105").Value = -3323.778

$ws.Range("H132").Value = 78723.625
$ws.Range("J132").Value = 78723.625
$ws.Range("L132").Value = 78723.625
$ws.Range("N132").Value = -88843.625

$ws.Range("H133").Value = 99249
$ws.Range("J133").Value = 99249
$ws.Range("L133").Value = 99249
$ws.Range("N133").Value = -109369

$ws.Range("H138").Value = 99764.336
$ws.Range("J138").Value = 99764.336
$ws.Range("L138").Value = 99764.336
$ws.Range("N138").Value = -110044.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2177.125
$ws.Range("I105").Value = 2235.4443
$ws.Range("K105").Value = 2235.4443
$ws.Range("M105").Value = -488.4443000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6445.6
$ws.Range("I137").Value = 5607.727
$ws.Range("J137").Value = 8749.75
$ws.Range("K137").Value = 16823.181
$ws.Range("L137").Value = 26249.25
$ws.Range("M137").Value = -11723.181
$ws.Range("N137").Value = -36449.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 30030
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H49").Value = 48500
$ws.Range("J49").Value = 48500
$ws.Range("L49").Value = 48500
$ws.Range("N49").Value = -48868

$ws.Range("H80").Value = 1200
$ws.Range("I80").Value = 1200
$ws.Range("K80").Value = 1200
$ws.Range("M80").Value = -202

$ws.Range("H83").Value = 1200
$ws.Range("I83").Value = 1200
$ws.Range("K83").Value = 6000
$ws.Range("M83").Value = -1008

$ws.Range("H122").Value = 1793
$ws.Range("I122").Value = 1798.6
$ws.Range("K122").Value = 5395.799999999999
$ws.Range("M122").Value = -2945.799999999999

$ws.Range("H126").Value = 22362.354
$ws.Range("J126").Value = 4499.75
$ws.Range("L126").Value = 13499.25
$ws.Range("N126").Value = -18439.25

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3314
$ws.Range("J22").Value = 3314
$ws.Range("L22").Value = 3314
$ws.Range("N22").Value = -3904

$ws.Range("H27").Value = 3314
$ws.Range("J27").Value = 3314
$ws.Range("L27").Value = 3314
$ws.Range("N27").Value = -3528

$ws.Range("H42").Value = 14770.3
$ws.Range("I42").Value = 12196
$ws.Range("J42").Value = 17344.6
$ws.Range("K42").Value = 12196
$ws.Range("L42").Value = 17344.6
$ws.Range("M42").Value = -11633
$ws.Range("N42").Value = -18470.6

$ws.Range("H49").Value = 14770.3
$ws.Range("I49").Value = 12196
$ws.Range("J49").Value = 17344.6
$ws.Range("K49").Value = 12196
$ws.Range("L49").Value = 17344.6
$ws.Range("M49").Value = -12049
$ws.Range("N49").Value = -17638.6

$ws.Range("H100").Value = 45378.793
$ws.Range("I100").Value = 59832.832
$ws.Range("J100").Value = 2016.6666
$ws.Range("K100").Value = 59832.832
$ws.Range("L100").Value = 2016.6666
$ws.Range("M100").Value = -59291.832
$ws.Range("N100").Value = -3098.6666

$ws.Range("H136").Value = 4178.2856
$ws.Range("I136").Value = 4039.158
$ws.Range("K136").Value = 12117.474
$ws.Range("M136").Value = -9567.474

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H126").Value = 4179.1177
$ws.Range("I126").Value = 4057.1428
$ws.Range("J126").Value = 4748.3335
$ws.Range("K126").Value = 12171.4284
$ws.Range("L126").Value = 14245.0005
$ws.Range("M126").Value = -9701.428400000001
$ws.Range("N126").Value = -19185.0005

$ws.Range("H132").Value = 2856.0454
$ws.Range("I132").Value = 2515.8572
$ws.Range("K132").Value = 7547.571599999999
$ws.Range("M132").Value = -5017.571599999999

$ws.Range("H133").Value = 89475
$ws.Range("J133").Value = 89475
$ws.Range("L133").Value = 89475
$ws.Range("N133").Value = -99595

$ws.Range("H139").Value = 52899.375
$ws.Range("J139").Value = 52899.375
$ws.Range("L139").Value = 52899.375
$ws.Range("N139").Value = -63179.375

$ws.Range("H140").Value = 59783.43
$ws.Range("J140").Value = 59783.43
$ws.Range("L140").Value = 59783.43
$ws.Range("N140").Value = -70143.42999999999

$ws.Range("H141").Value = 121917.8
$ws.Range("J141").Value = 121917.8
$ws.Range("L141").Value = 121917.8
$ws.Range("N141").Value = -132277.8
